$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet from "member" to "template_impor"
$ws.Name = "template_impor"

# 2. Add a new header column G: "is_outside_region", styled like the
#    existing headers (copy format from F1, the last header cell).
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "is_outside_region"

# 3. Column D (village_id) gets a date-like display format for most rows,
#    applied per contiguous block so each Copy/Paste keeps its own style.
#    (numFmt creation order matters so the custom numFmtIds land as
#    164=mm.yyyy, 165=yyyy-mm-dd, 166=m.yyyy)
$ws.Range("F1").Copy()
$ws.Range("D2:D12").PasteSpecial(-4122)
$ws.Range("D2:D12").NumberFormat = "mm.yyyy"

$ws.Range("F1").Copy()
$ws.Range("D16:D26").PasteSpecial(-4122)
$ws.Range("D16:D26").NumberFormat = "mm.yyyy"

# 4. Column F (event_date) gets a yyyy-mm-dd date format down through row 28.
$ws.Range("F1").Copy()
$ws.Range("F2:F28").PasteSpecial(-4122)
$ws.Range("F2:F28").NumberFormat = "yyyy-mm-dd"

# Rows 13, 28, 29 in column D use a plain thousands number format.
$ws.Range("F1").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").NumberFormat = "#,##0"

$ws.Range("F1").Copy()
$ws.Range("D28:D29").PasteSpecial(-4122)
$ws.Range("D28:D29").NumberFormat = "#,##0"

# Rows 14, 15, 30 in column D use a month.year format.
$ws.Range("F1").Copy()
$ws.Range("D14:D15").PasteSpecial(-4122)
$ws.Range("D14:D15").NumberFormat = "m.yyyy"

$ws.Range("F1").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").NumberFormat = "m.yyyy"

# 5. Row 26 gets a slightly taller custom row height.
$ws.Rows("26").RowHeight = 16.5
